# Update gh-pages to output generated at 456a3b4
#
# The two date-2024-02-07 events (row 2 and row 3) are removed from the
# "展览" and "全部类型" sheets; the remaining two events (originally rows 4
# and 5) move up to become rows 2 and 3, keeping their original index
# numbers (1 and 2) in column A, and their "想去人数" (column F) counts are
# bumped up (85 -> 87, 301 -> 302).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the old row 2 (2024-02-07 丽水·新年动漫狂欢盛典) and the old
    # row 3 (2024-02-07 龙泉·崩X铁X原ONLY). Deleting row 2 twice removes
    # both, shifting the remaining two data rows up.
    $ws.Rows(2).Delete() | Out-Null
    $ws.Rows(2).Delete() | Out-Null

    # Keep the original index numbering for the surviving rows.
    $ws.Range("A2").Value = 1
    $ws.Range("A3").Value = 2

    # Update "想去人数" counts for the surviving events.
    $ws.Range("F2").Value = 87
    $ws.Range("F3").Value = 302
}
